# Modify the bold "Planning and Identification." list heading so it reads
# "Planning and Identification Reports." — i.e. insert " Reports" just
# before the trailing period, keeping the existing bold/Arial run formatting.
#
# There are two paragraphs in the document that start with the phrase
# "Planning and Identification" — an earlier body paragraph that continues
# with more sentences ("...The first step in the process is...") and the
# short, bold numbered-list heading we actually need to change. We
# disambiguate by matching the paragraph whose *entire* text is exactly
# "Planning and Identification." and which is bold.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Planning and Identification.`r" -and $p.Range.Font.Bold) {
        $target = $p
    }
}

if ($target -ne $null) {
    $searchRange = $target.Range
    # Exclude the trailing paragraph mark so Find only ever matches inside
    # this one paragraph.
    $searchRange.MoveEnd(1, -1) | Out-Null

    $searchRange.Find.ClearFormatting()
    $searchRange.Find.Execute(
        "Planning and Identification.",  # FindText
        $false,                          # MatchCase
        $false,                          # MatchWholeWord
        $false,                          # MatchWildcards
        $false,                          # MatchSoundsLike
        $false,                          # MatchAllWordForms
        $true,                           # Forward
        1,                                # Wrap (wdFindContinue)
        $false,                          # Format
        "Planning and Identification Reports.",  # ReplaceWith
        2                                 # Replace (wdReplaceAll)
    ) | Out-Null
}
